$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 59.38757333333334
$ws.Range("H2").Value = 178.16272
$ws.Range("I2").Value = 0.1895016356445263
$ws.Range("J2").Value = 0.1895016356445263
$ws.Range("M2").Value = 227.11144
$ws.Range("N2").Value = 681.33432
$ws.Range("O2").Value = 0.8625743548356182
$ws.Range("P2").Value = 0.8625743548356182
$ws.Range("Q2").Value = 13487.59729783893
$ws.Range("R2").Value = 121388.3756805504
$ws.Range("S2").Value = 0.1634592511063717
$ws.Range("T2").Value = 0.1634592511063717

# Row 3
$ws.Range("G3").Value = 59.38757333333334
$ws.Range("H3").Value = 178.16272
$ws.Range("I3").Value = 0.1895016356445263
$ws.Range("J3").Value = 0.1895016356445263
$ws.Range("O3").Value = 0.001598666154760757
$ws.Range("P3").Value = 0.001598666154760757
$ws.Range("Q3").Value = 24.99745695918223
$ws.Range("R3").Value = 224.97711263264
$ws.Range("S3").Value = 0.0003029498511767088
$ws.Range("T3").Value = 0.0003029498511767088

# Row 4
$ws.Range("G4").Value = 59.38757333333334
$ws.Range("H4").Value = 178.16272
$ws.Range("I4").Value = 0.1895016356445263
$ws.Range("J4").Value = 0.1895016356445263
$ws.Range("M4").Value = 3.233093
$ws.Range("N4").Value = 9.699279000000001
$ws.Range("O4").Value = 0.01227935989749593
$ws.Range("P4").Value = 0.01227935989749593
$ws.Range("Q4").Value = 192.0055476309867
$ws.Range("R4").Value = 1728.04992867888
$ws.Range("S4").Value = 0.002326958785243282
$ws.Range("T4").Value = 0.002326958785243281

# Row 5
$ws.Range("G5").Value = 59.38757333333334
$ws.Range("H5").Value = 178.16272
$ws.Range("I5").Value = 0.1895016356445263
$ws.Range("J5").Value = 0.1895016356445263
$ws.Range("M5").Value = 32.52945966666667
$ws.Range("N5").Value = 97.588379
$ws.Range("O5").Value = 0.1235476191121251
$ws.Range("P5").Value = 0.1235476191121251
$ws.Range("Q5").Value = 1931.845671447876
$ws.Range("R5").Value = 17386.61104303088
$ws.Range("S5").Value = 0.02341247590173465
$ws.Range("T5").Value = 0.02341247590173465

# Row 6
$ws.Range("I6").Value = 0.6831820482914401
$ws.Range("J6").Value = 0.68318204829144
$ws.Range("M6").Value = 227.11144
$ws.Range("N6").Value = 681.33432
$ws.Range("O6").Value = 0.8625743548356182
$ws.Range("P6").Value = 0.8625743548356182
$ws.Range("Q6").Value = 48624.82752261064
$ws.Range("R6").Value = 437623.4477034958
$ws.Range("S6").Value = 0.5892953145402651
$ws.Range("T6").Value = 0.5892953145402651

# Row 7
$ws.Range("I7").Value = 0.6831820482914401
$ws.Range("J7").Value = 0.68318204829144
$ws.Range("O7").Value = 0.001598666154760757
$ws.Range("P7").Value = 0.001598666154760757
$ws.Range("S7").Value = 0.001092180018143654
$ws.Range("T7").Value = 0.001092180018143654

# Row 8
$ws.Range("I8").Value = 0.6831820482914401
$ws.Range("J8").Value = 0.68318204829144
$ws.Range("M8").Value = 3.233093
$ws.Range("N8").Value = 9.699279000000001
$ws.Range("O8").Value = 0.01227935989749593
$ws.Range("P8").Value = 0.01227935989749593
$ws.Range("Q8").Value = 692.2090295828331
$ws.Range("R8").Value = 6229.881266245497
$ws.Range("S8").Value = 0.008389038246479038
$ws.Range("T8").Value = 0.008389038246479037

# Row 9
$ws.Range("I9").Value = 0.6831820482914401
$ws.Range("J9").Value = 0.68318204829144
$ws.Range("M9").Value = 32.52945966666667
$ws.Range("N9").Value = 97.588379
$ws.Range("O9").Value = 0.1235476191121251
$ws.Range("P9").Value = 0.1235476191121251
$ws.Range("Q9").Value = 6964.5957319252
$ws.Range("R9").Value = 62681.3615873268
$ws.Range("S9").Value = 0.08440551548655233
$ws.Range("T9").Value = 0.08440551548655231

# Row 10
$ws.Range("G10").Value = 39.60717
$ws.Range("H10").Value = 118.82151
$ws.Range("I10").Value = 0.126383737825469
$ws.Range("J10").Value = 0.126383737825469
$ws.Range("M10").Value = 227.11144
$ws.Range("N10").Value = 681.33432
$ws.Range("O10").Value = 0.8625743548356182
$ws.Range("P10").Value = 0.8625743548356182
$ws.Range("Q10").Value = 8995.241413024802
$ws.Range("R10").Value = 80957.17271722322
$ws.Range("S10").Value = 0.1090153711165178
$ws.Range("T10").Value = 0.1090153711165178

# Row 11
$ws.Range("G11").Value = 39.60717
$ws.Range("H11").Value = 118.82151
$ws.Range("I11").Value = 0.126383737825469
$ws.Range("J11").Value = 0.126383737825469
$ws.Range("O11").Value = 0.001598666154760757
$ws.Range("P11").Value = 0.001598666154760757
$ws.Range("Q11").Value = 16.67147640118
$ws.Range("R11").Value = 150.04328761062
$ws.Range("S11").Value = 0.0002020454041737341
$ws.Range("T11").Value = 0.0002020454041737341

# Row 12
$ws.Range("G12").Value = 39.60717
$ws.Range("H12").Value = 118.82151
$ws.Range("I12").Value = 0.126383737825469
$ws.Range("J12").Value = 0.126383737825469
$ws.Range("M12").Value = 3.233093
$ws.Range("N12").Value = 9.699279000000001
$ws.Range("O12").Value = 0.01227935989749593
$ws.Range("P12").Value = 0.01227935989749593
$ws.Range("Q12").Value = 128.05366407681
$ws.Range("R12").Value = 1152.48297669129
$ws.Range("S12").Value = 0.001551911401949703
$ws.Range("T12").Value = 0.001551911401949703

# Row 13
$ws.Range("G13").Value = 39.60717
$ws.Range("H13").Value = 118.82151
$ws.Range("I13").Value = 0.126383737825469
$ws.Range("J13").Value = 0.126383737825469
$ws.Range("M13").Value = 32.52945966666667
$ws.Range("N13").Value = 97.588379
$ws.Range("O13").Value = 0.1235476191121251
$ws.Range("P13").Value = 0.1235476191121251
$ws.Range("Q13").Value = 1288.39983902581
$ws.Range("R13").Value = 11595.59855123229
$ws.Range("S13").Value = 0.01561440990282772
$ws.Range("T13").Value = 0.01561440990282772

# Row 14
$ws.Range("G14").Value = 0.292259
$ws.Range("H14").Value = 0.876777
$ws.Range("I14").Value = 0.0009325782385647279
$ws.Range("J14").Value = 0.000932578238564728
$ws.Range("M14").Value = 227.11144
$ws.Range("N14").Value = 681.33432
$ws.Range("O14").Value = 0.8625743548356182
$ws.Range("P14").Value = 0.8625743548356182
$ws.Range("Q14").Value = 66.37536234296
$ws.Range("R14").Value = 597.37826108664
$ws.Range("S14").Value = 0.0008044180724637073
$ws.Range("T14").Value = 0.0008044180724637074

# Row 15
$ws.Range("G15").Value = 0.292259
$ws.Range("H15").Value = 0.876777
$ws.Range("I15").Value = 0.0009325782385647279
$ws.Range("J15").Value = 0.000932578238564728
$ws.Range("O15").Value = 0.001598666154760757
$ws.Range("P15").Value = 0.001598666154760757
$ws.Range("Q15").Value = 0.1230178531193333
$ws.Range("R15").Value = 1.107160678074
$ws.Range("S15").Value = 0.000001490881266659833
$ws.Range("T15").Value = 0.000001490881266659833

# Row 16
$ws.Range("G16").Value = 0.292259
$ws.Range("H16").Value = 0.876777
$ws.Range("I16").Value = 0.0009325782385647279
$ws.Range("J16").Value = 0.000932578238564728
$ws.Range("M16").Value = 3.233093
$ws.Range("N16").Value = 9.699279000000001
$ws.Range("O16").Value = 0.01227935989749593
$ws.Range("P16").Value = 0.01227935989749593
$ws.Range("Q16").Value = 0.944900527087
$ws.Range("R16").Value = 8.504104743783001
$ws.Range("S16").Value = 0.00001145146382390911
$ws.Range("T16").Value = 0.00001145146382390911

# Row 17
$ws.Range("G17").Value = 0.292259
$ws.Range("H17").Value = 0.876777
$ws.Range("I17").Value = 0.0009325782385647279
$ws.Range("J17").Value = 0.000932578238564728
$ws.Range("M17").Value = 32.52945966666667
$ws.Range("N17").Value = 97.588379
$ws.Range("O17").Value = 0.1235476191121251
$ws.Range("P17").Value = 0.1235476191121251
$ws.Range("Q17").Value = 9.507027352720334
$ws.Range("R17").Value = 85.56324617448301
$ws.Range("S17").Value = 0.0001152178210104516
$ws.Range("T17").Value = 0.0001152178210104516
